$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.841.81'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '3.833.37'
$ws.Range("E3").Value = '  +3.84%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '410.77'
$ws.Range("E5").Value = '  -2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.67'
$ws.Range("E6").Value = '  +1.35%  '
$ws.Range("D7").Value = '3.820.64'
$ws.Range("E7").Value = '  +3.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.615'
$ws.Range("E8").Value = '  -4.38%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.734'
$ws.Range("E10").Value = '  -4.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.169'
$ws.Range("E11").Value = '  -6.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000369'
$ws.Range("E12").Value = '  -6.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.92'
$ws.Range("E13").Value = '  -5.05%  '
$ws.Range("D14").Value = '4.469.95'
$ws.Range("E14").Value = '  +4.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '9.97'
$ws.Range("E15").Value = '  -6.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.11'
$ws.Range("E16").Value = '  +14.00%  '
$ws.Range("D17").Value = '3.867.74'
$ws.Range("E17").Value = '  +4.78%  '
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19.60'
$ws.Range("E19").Value = '  -4.72%  '
$ws.Range("D20").Value = '67.392.16'
$ws.Range("E20").Value = '  +1.13%  '
$ws.Range("E21").Value = '  -5.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '412.62'
$ws.Range("E22").Value = '  -7.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.44'
$ws.Range("E23").Value = '  -12.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.47'
$ws.Range("E24").Value = '  -4.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.06'
$ws.Range("E25").Value = '  -2.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '36.64'
$ws.Range("E26").Value = '  -1.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.75'
$ws.Range("E27").Value = '  +13.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.12'
$ws.Range("E28").Value = '  -6.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.45'
$ws.Range("E29").Value = '  -7.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '681.38'
$ws.Range("E30").Value = '  +4.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.46'
$ws.Range("E31").Value = '  -2.10%  '
$ws.Range("E32").Value = '  -2.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.73'
$ws.Range("E33").Value = '  -1.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.20'
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.152'
$ws.Range("E35").Value = '  -9.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.83'
$ws.Range("E36").Value = '  -6.30%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0792'
$ws.Range("E38").Value = '  +7.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '54.94'
$ws.Range("E39").Value = '  -4.14%  '
$ws.Range("B40").Value = 'ThetaToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.09'
$ws.Range("E40").Value = '  +0.75%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0458'
$ws.Range("E41").Value = '  -7.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.137'
$ws.Range("E43").Value = '  -8.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '148.96'
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.47'
$ws.Range("E45").Value = '  +3.04%  '
$ws.Range("B46").Value = 'LidoDAOToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.33'
$ws.Range("E46").Value = '  -2.61%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.17'
$ws.Range("E47").Value = '  +17.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.87'
$ws.Range("E48").Value = '  -8.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.07'
$ws.Range("E49").Value = '  -1.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.80'
$ws.Range("E50").Value = '  -3.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.56'
$ws.Range("E51").Value = '  -4.24%  '
